$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "1.000", "93.60") are stored as text, not re-parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.699.43"
$ws.Range("D3").Value = "1.775.28"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "326.16"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.4595"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("D8").Value = "0.3589"
$ws.Range("D9").Value = "0.07479"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "41.95"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "1.102"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "20.83"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "6.036"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "7.224"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "1.776.50"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "93.60"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "0.00001058"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "0.06413"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "17.07"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "5.790"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "27.781.83"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").Value = "2.080"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "164.59"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").Value = "20.36"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "1.980.89"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "2.163"
$ws.Range("E29").Value = "  +4.75%  "
$ws.Range("D30").Value = "126.22"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "1.097"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").Value = "0.09222"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "3.679"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "5.533"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "11.80"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "0.02292"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "0.06108"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").Value = "0.2089"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "0.6308"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "4.962"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "1.394"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "7.764"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "13.17"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "3.729"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "0.5884"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "122.29"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").Value = "1.949"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "0.06943"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "1.138"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "72.30"
$ws.Range("E51").Value = "  +0.91%  "

# Restore the original (default/Normal) style on column D so no stray
# number-format style is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
